$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.298.28"
$ws.Range("E2").Value = "  +0.27%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.859.68"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - XRP
$ws.Range("D5").Value = "0.7031"
$ws.Range("E5").Value = "  +0.45%  "

# Row 6 - BNB
$ws.Range("D6").Value = "238.29"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7 - USDC
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.07899"
$ws.Range("E8").Value = "  +2.97%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.3043"
$ws.Range("E9").Value = "  -0.12%  "

# Row 10 - Solana
$ws.Range("D10").Value = "24.41"
$ws.Range("E10").Value = "  +4.78%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.08181"
$ws.Range("E11").Value = "  +0.36%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.863.82"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "0.7229"
$ws.Range("E13").Value = "  +0.87%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "5.217"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15 - Litecoin (D unchanged)
$ws.Range("E15").Value = "  +0.29%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "29.314.25"
$ws.Range("E16").Value = "  +0.27%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "5.816"
$ws.Range("E17").Value = "  +1.12%  "

# Row 18 - ShibaInu (E unchanged)
$ws.Range("D18").Value = "0.000007820"

# Row 19 - Avalanche (D unchanged)
$ws.Range("E19").Value = "  +0.06%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "238.41"
$ws.Range("E20").Value = "  +0.40%  "

# Rows 21/22 swap: Dai now at row 21, WrappedliquidstakedEther2.0 now at row 22
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.114.19"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23 - BinanceUSD (D unchanged)
$ws.Range("E23").Value = "  -0.05%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "7.565"
$ws.Range("E24").Value = "  +1.49%  "

# Row 25 - Monero
$ws.Range("D25").Value = "162.05"
$ws.Range("E25").Value = "  -0.13%  "

# Row 26 - Cosmos (D unchanged)
$ws.Range("E26").Value = "  -1.11%  "

# Row 27 - Stellar
$ws.Range("D27").Value = "0.1431"
$ws.Range("E27").Value = "  -2.82%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "18.13"
$ws.Range("E28").Value = "  +0.59%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "1.914"
$ws.Range("E29").Value = "  -4.55%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "1.385"
$ws.Range("E30").Value = "  -2.50%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.476"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.326"
$ws.Range("E32").Value = "  -2.38%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "4.056"
$ws.Range("E33").Value = "  +1.24%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.05185"
$ws.Range("E34").Value = "  -0.23%  "

# Row 35 - ARBITRUM (D unchanged)
$ws.Range("E35").Value = "  +0.83%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.7140"
$ws.Range("E36").Value = "  +0.55%  "

# Row 37 - Frax
$ws.Range("D37").Value = "0.9988"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 - HuobiToken
$ws.Range("D38").Value = "2.676"
$ws.Range("E38").Value = "  +0.69%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.01852"
$ws.Range("E39").Value = "  +0.01%  "

# Row 40 - MXToken (D unchanged)
$ws.Range("E40").Value = "  -1.31%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.155.90"
$ws.Range("E41").Value = "  +0.60%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.9216"
$ws.Range("E42").Value = "  -1.23%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "5.960"
$ws.Range("E43").Value = "  +1.73%  "

# Row 44 - TheSandbox
$ws.Range("D44").Value = "0.4258"
$ws.Range("E44").Value = "  -0.59%  "

# Row 45 - Aave
$ws.Range("D45").Value = "70.91"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46 - PaxDollar (D unchanged)
$ws.Range("E46").Value = "  -0.04%  "

# Row 47 - Quant
$ws.Range("D47").Value = "101.81"
$ws.Range("E47").Value = "  -1.43%  "

# Row 48 - Mantle (D unchanged)
$ws.Range("E48").Value = "  -2.95%  "

# Row 49 - RenderToken (D unchanged)
$ws.Range("E49").Value = "  -2.52%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "9.188"
$ws.Range("E50").Value = "  +0.43%  "

# Row 51 - Aptos
$ws.Range("D51").Value = "7.001"
$ws.Range("E51").Value = "  +0.56%  "
